$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($RangeAddr, $Val)
    $c = $ws.Range($RangeAddr)
    $escaped = $Val.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" "69.106.49"
$ws.Range("E2").Value = "  +0.35%  "
Set-TextValue "D3" "3.785.57"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "602.47"
$ws.Range("E5").Value = "  -0.03%  "
Set-TextValue "D6" "164.20"
$ws.Range("E6").Value = "  -2.84%  "
Set-TextValue "D7" "3.781.67"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.537"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  +2.86%  "
Set-TextValue "D11" "6.31"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -0.31%  "
Set-TextValue "D13" "37.34"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -0.79%  "
Set-TextValue "D15" "4.418.88"
$ws.Range("E15").Value = "  +1.08%  "
Set-TextValue "D16" "3.791.95"
$ws.Range("E16").Value = "  +1.29%  "
Set-TextValue "D17" "69.218.76"
$ws.Range("E18").Value = "  +2.04%  "
Set-TextValue "D19" "17.36"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("E21").Value = "  +3.23%  "
Set-TextValue "D22" "491.94"
$ws.Range("E22").Value = "  -0.29%  "
Set-TextValue "D23" "0.722"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -2.23%  "
Set-TextValue "D25" "84.66"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  -3.08%  "
Set-TextValue "D27" "12.23"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.10%  "
Set-TextValue "D31" "8.12"
$ws.Range("E31").Value = "  +2.17%  "
Set-TextValue "D32" "2.40"
$ws.Range("E32").Value = "  -4.84%  "
Set-TextValue "D33" "3.935.63"
$ws.Range("E33").Value = "  +1.14%  "
Set-TextValue "D34" "31.85"
$ws.Range("E34").Value = "  +0.26%  "
Set-TextValue "D35" "3.732.06"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("E38").Value = "  +0.48%  "
Set-TextValue "D39" "5.93"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.22%  "
Set-TextValue "D42" "3.06"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D43" "425.48"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D44" "48.41"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D46" "1.00"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D47" "8.41"
$ws.Range("E47").Value = "  -0.62%  "
Set-TextValue "D48" "142.43"
$ws.Range("E48").Value = "  +0.61%  "
Set-TextValue "D49" "2.818.06"
$ws.Range("E49").Value = "  +1.38%  "
Set-TextValue "D50" "39.90"
$ws.Range("E50").Value = "  -1.65%  "
Set-TextValue "D51" "1.30"
$ws.Range("E51").Value = "  +7.76%  "

$excel.CutCopyMode = $false
